# Weekly refresh of the Fruta/Hortaliza price sheet: the data for rows 2-11
# (Albahaca @ Vega Monumental Concepción) is re-shuffled across the date range
# -- existing per-date records simply move to a different row while columns
# A,B,C,E,F,G,H,I,N,Q,R stay constant (same market/category for every row).
# Row 4 is unaffected. Only D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen) and
# P (Precio $/Kg) move between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg)

    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
}

# New row 2 <= old row 6
Set-Row 2  44671 150 3500 4000 3733 "Región Metropolitana"   622
# New row 3 <= old row 9
Set-Row 3  44630  90 2500 3000 2722 "Región Metropolitana"   454
# Row 4 unchanged
Set-Row 4  44672 140 3000 3500 3286 "Región Metropolitana"   548
# New row 5 <= old row 8
Set-Row 5  44658 180 2500 3000 2778 "Región Metropolitana"   463
# New row 6 <= old row 2
Set-Row 6  44631 110 3000 3500 3273 "Provincia de Chacabuco" 546
# New row 7 <= old row 5
Set-Row 7  44637 170 2800 3000 2906 "Región Metropolitana"   484
# New row 8 <= old row 10
Set-Row 8  44644 140 2500 3000 2786 "Provincia de Chacabuco" 464
# New row 9 <= old row 7
Set-Row 9  44643  90 2800 3000 2911 "Región Metropolitana"   485
# New row 10 <= old row 11
Set-Row 10 44650 130 3000 3500 3308 "Región Metropolitana"   551
# New row 11 <= old row 3
Set-Row 11 44659  90 2500 3000 2722 "Región Metropolitana"   454

Write-Output "Weekly reshuffle applied"
